$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb1"
$ws.Cells.Item(2, 3).Value = "Epha4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 7.079689000000001
$ws.Cells.Item(2, 8).Value = 21.239067
$ws.Cells.Item(2, 9).Value = 0.5033576067109902
$ws.Cells.Item(2, 10).Value = 0.5033576067109902
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 5.685057
$ws.Cells.Item(2, 14).Value = 17.055171
$ws.Cells.Item(2, 15).Value = 0.3604606774420115
$ws.Cells.Item(2, 16).Value = 0.3604606774420115
$ws.Cells.Item(2, 17).Value = 40.24843550727301
$ws.Cells.Item(2, 18).Value = 362.235919565457
$ws.Cells.Item(2, 19).Value = 0.1814406239106331
$ws.Cells.Item(2, 20).Value = 0.1814406239106331

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb1"
$ws.Cells.Item(3, 3).Value = "Epha4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 7.079689000000001
$ws.Cells.Item(3, 8).Value = 21.239067
$ws.Cells.Item(3, 9).Value = 0.5033576067109902
$ws.Cells.Item(3, 10).Value = 0.5033576067109902
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.775186333333332
$ws.Cells.Item(3, 14).Value = 26.325559
$ws.Cells.Item(3, 15).Value = 0.556390131249909
$ws.Cells.Item(3, 16).Value = 0.5563901312499091
$ws.Cells.Item(3, 17).Value = 62.12559015705033
$ws.Cells.Item(3, 18).Value = 559.130311413453
$ws.Cells.Item(3, 19).Value = 0.2800632048635679
$ws.Cells.Item(3, 20).Value = 0.2800632048635679

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb1"
$ws.Cells.Item(4, 3).Value = "Epha4"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 7.079689000000001
$ws.Cells.Item(4, 8).Value = 21.239067
$ws.Cells.Item(4, 9).Value = 0.5033576067109902
$ws.Cells.Item(4, 10).Value = 0.5033576067109902
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.04169666666666667
$ws.Cells.Item(4, 14).Value = 0.12509
$ws.Cells.Item(4, 15).Value = 0.002643774497553922
$ws.Cells.Item(4, 16).Value = 0.002643774497553922
$ws.Cells.Item(4, 17).Value = 0.2951994323366667
$ws.Cells.Item(4, 18).Value = 2.656794891030001
$ws.Cells.Item(4, 19).Value = 0.001330764003772293
$ws.Cells.Item(4, 20).Value = 0.001330764003772293

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efnb1"
$ws.Cells.Item(5, 3).Value = "Epha4"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.079689000000001
$ws.Cells.Item(5, 8).Value = 21.239067
$ws.Cells.Item(5, 9).Value = 0.5033576067109902
$ws.Cells.Item(5, 10).Value = 0.5033576067109902
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.269702666666667
$ws.Cells.Item(5, 14).Value = 3.809108
$ws.Cells.Item(5, 15).Value = 0.08050541681052542
$ws.Cells.Item(5, 16).Value = 0.08050541681052542
$ws.Cells.Item(5, 17).Value = 8.989100002470668
$ws.Cells.Item(5, 18).Value = 80.901900022236
$ws.Cells.Item(5, 19).Value = 0.04052301393301679
$ws.Cells.Item(5, 20).Value = 0.04052301393301679

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb1"
$ws.Cells.Item(6, 3).Value = "Epha4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.058683666666667
$ws.Cells.Item(6, 8).Value = 12.176051
$ws.Cells.Item(6, 9).Value = 0.2885676612136944
$ws.Cells.Item(6, 10).Value = 0.2885676612136945
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 5.685057
$ws.Cells.Item(6, 14).Value = 17.055171
$ws.Cells.Item(6, 15).Value = 0.3604606774420115
$ws.Cells.Item(6, 16).Value = 0.3604606774420115
$ws.Cells.Item(6, 17).Value = 23.073847989969
$ws.Cells.Item(6, 18).Value = 207.664631909721
$ws.Cells.Item(6, 19).Value = 0.1040172946489452
$ws.Cells.Item(6, 20).Value = 0.1040172946489452

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb1"
$ws.Cells.Item(7, 3).Value = "Epha4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.058683666666667
$ws.Cells.Item(7, 8).Value = 12.176051
$ws.Cells.Item(7, 9).Value = 0.2885676612136944
$ws.Cells.Item(7, 10).Value = 0.2885676612136945
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 8.775186333333332
$ws.Cells.Item(7, 14).Value = 26.325559
$ws.Cells.Item(7, 15).Value = 0.556390131249909
$ws.Cells.Item(7, 16).Value = 0.5563901312499091
$ws.Cells.Item(7, 17).Value = 35.61570544305656
$ws.Cells.Item(7, 18).Value = 320.541348987509
$ws.Cells.Item(7, 19).Value = 0.1605561988971667
$ws.Cells.Item(7, 20).Value = 0.1605561988971668

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efnb1"
$ws.Cells.Item(8, 3).Value = "Epha4"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.058683666666667
$ws.Cells.Item(8, 8).Value = 12.176051
$ws.Cells.Item(8, 9).Value = 0.2885676612136944
$ws.Cells.Item(8, 10).Value = 0.2885676612136945
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04169666666666667
$ws.Cells.Item(8, 14).Value = 0.12509
$ws.Cells.Item(8, 15).Value = 0.002643774497553922
$ws.Cells.Item(8, 16).Value = 0.002643774497553922
$ws.Cells.Item(8, 17).Value = 0.1692335799544445
$ws.Cells.Item(8, 18).Value = 1.52310221959
$ws.Cells.Item(8, 19).Value = 0.0007629078235355454
$ws.Cells.Item(8, 20).Value = 0.0007629078235355456

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efnb1"
$ws.Cells.Item(9, 3).Value = "Epha4"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.058683666666667
$ws.Cells.Item(9, 8).Value = 12.176051
$ws.Cells.Item(9, 9).Value = 0.2885676612136944
$ws.Cells.Item(9, 10).Value = 0.2885676612136945
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.269702666666667
$ws.Cells.Item(9, 14).Value = 3.809108
$ws.Cells.Item(9, 15).Value = 0.08050541681052542
$ws.Cells.Item(9, 16).Value = 0.08050541681052542
$ws.Cells.Item(9, 17).Value = 5.153321474723112
$ws.Cells.Item(9, 18).Value = 46.37989327250801
$ws.Cells.Item(9, 19).Value = 0.02323125984404696
$ws.Cells.Item(9, 20).Value = 0.02323125984404696

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Efnb1"
$ws.Cells.Item(10, 3).Value = "Epha4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.463315
$ws.Cells.Item(10, 8).Value = 1.389945
$ws.Cells.Item(10, 9).Value = 0.03294115455541936
$ws.Cells.Item(10, 10).Value = 0.03294115455541936
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 5.685057
$ws.Cells.Item(10, 14).Value = 17.055171
$ws.Cells.Item(10, 15).Value = 0.3604606774420115
$ws.Cells.Item(10, 16).Value = 0.3604606774420115
$ws.Cells.Item(10, 17).Value = 2.633972183955
$ws.Cells.Item(10, 18).Value = 23.705749655595
$ws.Cells.Item(10, 19).Value = 0.01187399088676847
$ws.Cells.Item(10, 20).Value = 0.01187399088676847

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Efnb1"
$ws.Cells.Item(11, 3).Value = "Epha4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.463315
$ws.Cells.Item(11, 8).Value = 1.389945
$ws.Cells.Item(11, 9).Value = 0.03294115455541936
$ws.Cells.Item(11, 10).Value = 0.03294115455541936
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 8.775186333333332
$ws.Cells.Item(11, 14).Value = 26.325559
$ws.Cells.Item(11, 15).Value = 0.556390131249909
$ws.Cells.Item(11, 16).Value = 0.5563901312499091
$ws.Cells.Item(11, 17).Value = 4.065675456028333
$ws.Cells.Item(11, 18).Value = 36.59107910425499
$ws.Cells.Item(11, 19).Value = 0.01832813330661332
$ws.Cells.Item(11, 20).Value = 0.01832813330661332

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Efnb1"
$ws.Cells.Item(12, 3).Value = "Epha4"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.463315
$ws.Cells.Item(12, 8).Value = 1.389945
$ws.Cells.Item(12, 9).Value = 0.03294115455541936
$ws.Cells.Item(12, 10).Value = 0.03294115455541936
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.04169666666666667
$ws.Cells.Item(12, 14).Value = 0.12509
$ws.Cells.Item(12, 15).Value = 0.002643774497553922
$ws.Cells.Item(12, 16).Value = 0.002643774497553922
$ws.Cells.Item(12, 17).Value = 0.01931869111666666
$ws.Cells.Item(12, 18).Value = 0.17386822005
$ws.Cells.Item(12, 19).Value = [double]"8.708898433359992e-05"
$ws.Cells.Item(12, 20).Value = [double]"8.708898433359993e-05"

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Efnb1"
$ws.Cells.Item(13, 3).Value = "Epha4"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.463315
$ws.Cells.Item(13, 8).Value = 1.389945
$ws.Cells.Item(13, 9).Value = 0.03294115455541936
$ws.Cells.Item(13, 10).Value = 0.03294115455541936
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 1.269702666666667
$ws.Cells.Item(13, 14).Value = 3.809108
$ws.Cells.Item(13, 15).Value = 0.08050541681052542
$ws.Cells.Item(13, 16).Value = 0.08050541681052542
$ws.Cells.Item(13, 17).Value = 0.5882722910066667
$ws.Cells.Item(13, 18).Value = 5.29445061906
$ws.Cells.Item(13, 19).Value = 0.002651941377703974
$ws.Cells.Item(13, 20).Value = 0.002651941377703974

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Efnb1"
$ws.Cells.Item(14, 3).Value = "Epha4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.463241333333333
$ws.Cells.Item(14, 8).Value = 7.389724
$ws.Cells.Item(14, 9).Value = 0.175133577519896
$ws.Cells.Item(14, 10).Value = 0.175133577519896
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 5.685057
$ws.Cells.Item(14, 14).Value = 17.055171
$ws.Cells.Item(14, 15).Value = 0.3604606774420115
$ws.Cells.Item(14, 16).Value = 0.3604606774420115
$ws.Cells.Item(14, 17).Value = 14.003667384756
$ws.Cells.Item(14, 18).Value = 126.033006462804
$ws.Cells.Item(14, 19).Value = 0.06312876799566475
$ws.Cells.Item(14, 20).Value = 0.06312876799566475

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Efnb1"
$ws.Cells.Item(15, 3).Value = "Epha4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.463241333333333
$ws.Cells.Item(15, 8).Value = 7.389724
$ws.Cells.Item(15, 9).Value = 0.175133577519896
$ws.Cells.Item(15, 10).Value = 0.175133577519896
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 8.775186333333332
$ws.Cells.Item(15, 14).Value = 26.325559
$ws.Cells.Item(15, 15).Value = 0.556390131249909
$ws.Cells.Item(15, 16).Value = 0.5563901312499091
$ws.Cells.Item(15, 17).Value = 21.61540168396844
$ws.Cells.Item(15, 18).Value = 194.538615155716
$ws.Cells.Item(15, 19).Value = 0.09744259418256103
$ws.Cells.Item(15, 20).Value = 0.09744259418256106

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Efnb1"
$ws.Cells.Item(16, 3).Value = "Epha4"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.463241333333333
$ws.Cells.Item(16, 8).Value = 7.389724
$ws.Cells.Item(16, 9).Value = 0.175133577519896
$ws.Cells.Item(16, 10).Value = 0.175133577519896
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.04169666666666667
$ws.Cells.Item(16, 14).Value = 0.12509
$ws.Cells.Item(16, 15).Value = 0.002643774497553922
$ws.Cells.Item(16, 16).Value = 0.002643774497553922
$ws.Cells.Item(16, 17).Value = 0.1027089527955556
$ws.Cells.Item(16, 18).Value = 0.92438057516
$ws.Cells.Item(16, 19).Value = 0.0004630136859124838
$ws.Cells.Item(16, 20).Value = 0.0004630136859124839

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Efnb1"
$ws.Cells.Item(17, 3).Value = "Epha4"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2.463241333333333
$ws.Cells.Item(17, 8).Value = 7.389724
$ws.Cells.Item(17, 9).Value = 0.175133577519896
$ws.Cells.Item(17, 10).Value = 0.175133577519896
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 1.269702666666667
$ws.Cells.Item(17, 14).Value = 3.809108
$ws.Cells.Item(17, 15).Value = 0.08050541681052542
$ws.Cells.Item(17, 16).Value = 0.08050541681052542
$ws.Cells.Item(17, 17).Value = 3.127584089576889
$ws.Cells.Item(17, 18).Value = 28.148256806192
$ws.Cells.Item(17, 19).Value = 0.01409920165575769
$ws.Cells.Item(17, 20).Value = 0.01409920165575769
